# Sheet1 currently holds a 4-column table (EmpName, City, EmpID, EmailID).
# The edit rotates the first three columns one position to the right
# (old col A -> new col B, old col B -> new col C, old col C -> new col A),
# so the table becomes (EmpID, EmpName, City, EmailID). Column D (EmailID,
# with its hyperlinks) is untouched.
#
# Doing this as a real column Cut + Insert (rather than copying cell values)
# makes Excel carry the column widths and cell formatting along with the
# data, exactly like a user dragging the "EmpID" column to the front.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cut column C (EmpID) and insert it before column A, shifting
# A (EmpName) -> B and B (City) -> C.
$ws.Columns.Item(3).Cut()
$ws.Columns.Item(1).Insert()

# Move the active selection to D14, matching the saved cursor position.
$ws.Range("D14").Select() | Out-Null
